$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (CN3): description text updated to the new connector description.
# Leading "'" forces a text entry (matches the quote-prefixed text style the
# other populated cells in this column already use) without the apostrophe
# itself ending up in the stored value.
$ws.Range("E5").Value = "'Female Header Micro Match 1,27 mm Wire to Board Connector Series 369"

# Row 6 (CN4): fill in manufacturer / part number / description that were
# previously blank.
$ws.Range("C6").Value = "'MPE-Garry"
$ws.Range("D6").Value = "369-1-024-0-NTX-KT0"
$ws.Range("E6").Value = "'Female Header Micro Match 1,27 mm Wire to Board Connector Series 369"

# The new, longer description wraps to two lines, so the row grows to match
# the other two-line rows in the sheet.
$ws.Rows.Item(6).RowHeight = 29

$ws.Range("E6").Select()
